# The original row 13 ("1720367 - Teresa Cristina Brazil de Paiva", under the
# "Docentes responsáveis:" label) is removed from the sheet, which shifts every
# row below it up by one. On top of that shift, four of the now-relocated rows
# get their data-column (B/C) text replaced with new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Docentes responsáveis:" value row; everything below moves up.
$ws.Rows(13).Delete()

# Row 10 ("Objetivos:") now shows the professor string that used to live in the
# deleted row.
$ws.Range("B10:C10").Value = "1720367 - Teresa Cristina Brazil de Paiva"

# Row 13 ("Programa resumido:", formerly row 14) now just says "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 ("Programa:", formerly row 16) now shows the activation date.
$ws.Range("B15:C15").Value = "01/01/2018"

# Row 18 ("Método:", formerly row 19) now shows the professor string too.
$ws.Range("B18:C18").Value = "1720367 - Teresa Cristina Brazil de Paiva"
